# Day 11 code uploaded - add "Task 3" date column and mark attendance ("p")
# for Week 3 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 3")

# Add the new date header for Task 3 (F4) - 21 June 2024 (serial 45464)
$ws.Range("F4").Value = 45464

# Mark attendance ("p") for the students that submitted Task 3
$rows = @(5, 6, 8, 9, 10, 15, 16, 17, 18, 19, 21, 23, 25)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "p"
}

# Update the active selection to reflect where the editor left off
$ws.Range("J19").Select()
